# Auto-generated edit script applying the cryptos.xlsx data refresh
# (GitHub Actions scheduled update: prices + 1h volume deltas,
# plus a few rank swaps where two coins traded places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.900.98"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.889.35"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7735"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.88"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3114"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.62"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07172"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08600"
$ws.Range("E11").Value = "  +6.35%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7629"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.940.24"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.373"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.77"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.179"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "29.951.03"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.77"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.36"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007829"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "2.205.58"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.005"
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1644"
$ws.Range("E25").Value = "  +4.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.382"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.05"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.76"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.039"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.443"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.536"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.100"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05432"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.239"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7457"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.694"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.780"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.115.67"
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4468"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.081"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.09"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8516"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.53"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.865"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.633"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("D50").Value = "2.098.28"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.981"
$ws.Range("E51").Value = "  -0.74%  "
